$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (columns B:G) with corrected values ---
$ws.Range("B2").Value = 0.5224929946896146
$ws.Range("C2").Value = 1.898968271571221
$ws.Range("D2").Value = 15.67148254258704
$ws.Range("E2").Value = 3.958722337142003
$ws.Range("F2").Value = 3.967451592372942
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = 0.4961555363655573
$ws.Range("C3").Value = 1.850890984938442
$ws.Range("D3").Value = 15.31652666909897
$ws.Range("E3").Value = 3.913633435708941
$ws.Range("F3").Value = 3.925922114836308
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = 0.6562009383051112
$ws.Range("C4").Value = 1.755547849851764
$ws.Range("D4").Value = 14.91514394174393
$ws.Range("E4").Value = 3.862012939095871
$ws.Range("F4").Value = 3.849856378757894
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = 0.5233911015648834
$ws.Range("C5").Value = 1.855814880858909
$ws.Range("D5").Value = 15.64356210993208
$ws.Range("E5").Value = 3.955194320122853
$ws.Range("F5").Value = 3.966808199306541
$ws.Range("G5").Value = 43

# --- Append new rows 6-11 (Q4..Q9) ---
$newRows = @(
    @{ Row=6;  Label="Q4"; B=0.7469272014119552; C=1.819658143763748; D=15.50982749350374; E=3.938251832159003; F=3.913643877732932; G=42 },
    @{ Row=7;  Label="Q5"; B=0.6021480811658927; C=1.875207305608213; D=16.24141422234426; E=4.030063798793297; F=4.034328036793712; G=41 },
    @{ Row=8;  Label="Q6"; B=0.7822106663771212; C=1.836403048259028; D=16.36623743788394; E=4.045520663386103; F=4.019744080550121; G=40 },
    @{ Row=9;  Label="Q7"; B=0.5940989821051248; C=1.932689778747934; D=16.9873845567875;  E=4.121575494490851; F=4.131849360674751; G=39 },
    @{ Row=10; Label="Q8"; B=0.830600775663463;  C=1.865149275164063; D=17.27800395215996; E=4.156681843990463; F=4.12752147190327;  G=38 },
    @{ Row=11; Label="Q9"; B=0.5950944418594691; C=1.827659872411423; D=17.20024814321981; E=4.147318186879301; F=4.161016494713655; G=37 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Bring over the label cell formatting (border/bold/alignment) from the row above
    $ws.Range("A" + ($row - 1)).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)

    $ws.Range("A" + $row).Value = $r.Label
    $ws.Range("B" + $row).Value = $r.B
    $ws.Range("C" + $row).Value = $r.C
    $ws.Range("D" + $row).Value = $r.D
    $ws.Range("E" + $row).Value = $r.E
    $ws.Range("F" + $row).Value = $r.F
    $ws.Range("G" + $row).Value = $r.G
}

$excel.CutCopyMode = $false
